$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the non-hyperlink data cells (columns B and C) ---
$ws.Range("B2").Value = "admin"
$ws.Range("C2").Value = "Pass"

$ws.Range("B3").Value = "admin"
$ws.Range("C3").Value = "Fail"

$ws.Range("B4").Value = "admin1"
$ws.Range("C4").Value = "Fail"

$ws.Range("B5").Value = "admin1"
$ws.Range("C5").Value = "Fail"

# --- Column A values + hyperlinks (reproduces Excel's auto-hyperlink /
#     copy-paste leftover pattern captured in the target workbook) ---

# A2: admin@yourstore.com, with its own hyperlink (rId1)
$ws.Range("A2").Value = "admin@yourstore.com"
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:admin@yourstore.com")

# Leftover wide hyperlink ranging A3:A5 (rId2), carrying a stale display
# value of "admin@yourstore.com" (this is what is left on disk in the
# source workbook from an earlier paste of A2 over A3:A5).
$ws.Hyperlinks.Add($ws.Range("A3:A5"), "mailto:admin@yourstore.com", "", "", "admin@yourstore.com")

# A3's real text is admin1@yourstore.com -- update it after the wide
# hyperlink so the stale display= text on rId2 is untouched.
$ws.Range("A3").Value = "admin1@yourstore.com"

# A4 only gets the hyperlink-style formatting from the wide rId2 range
# above (no dedicated relationship of its own in the target), but its
# text is admin@yourstore.com.
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:admin@yourstore.com")
$ws.Hyperlinks.Item($ws.Hyperlinks.Count).Delete()
$ws.Range("A4").Value = "admin@yourstore.com"

# A3 gets its own dedicated hyperlink (rId3) on top of the wide one.
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:admin1@yourstore.com")

# A5 gets its own dedicated hyperlink (rId4).
$ws.Range("A5").Value = "admin1@yourstore.com"
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:admin1@yourstore.com")

# --- Selection moved to B11 in the saved view state ---
$ws.Range("B11").Select()
